$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are stored as text (number-looking strings) in the sheet, so force
# the cells to keep a text number format before writing the new values -
# this prevents Excel from auto-converting them into numeric cells.
$cells = @("C3", "D3", "C4", "D4", "E4", "C5", "D5", "C6", "D6", "E6")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C3").Value = "19"
$ws.Range("D3").Value = "17"

$ws.Range("C4").Value = "21"
$ws.Range("D4").Value = "26"
$ws.Range("E4").Value = "4"

$ws.Range("C5").Value = "15"
$ws.Range("D5").Value = "24"

$ws.Range("C6").Value = "14"
$ws.Range("D6").Value = "11"
$ws.Range("E6").Value = "2"
